$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 592 - shifts existing rows 592:617 down to 593:618
# and extends the sheet's used range to A1:R618.
$ws.Rows(592).Insert()

# Populate the newly inserted row 592 with a new weekly price record
# (same market/category/quality as the surrounding rows).
$ws.Range("A592").Value = 3
$ws.Range("B592").Value = 'Femacal de La Calera'
$ws.Range("C592").Value = 'Coquimbo'
$ws.Range("D592").Value = 45147
$ws.Range("E592").Value = 5
$ws.Range("F592").Value = 100112043
$ws.Range("G592").Value = 'Pepino ensalada'
$ws.Range("H592").Value = 'Sin especificar'
$ws.Range("I592").Value = 'Primera'
$ws.Range("J592").Value = 90
$ws.Range("K592").Value = 10000
$ws.Range("L592").Value = 11000
$ws.Range("M592").Value = 10389
$ws.Range("N592").Value = '$/caja 60 unidades'
$ws.Range("O592").Value = 'Región de Arica y Parinacota'
$ws.Range("P592").Value = 173
$ws.Range("Q592").Value = 60
$ws.Range("R592").Value = 'Hortaliza'
